$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.546.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.315.40"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.335.31"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.728.17"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.612.34"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.313.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.27%  "
$ws.Range("E28").Value = "  +12.75%  "
$ws.Range("E29").Value = "  +5.54%  "
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.65"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.40"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.921"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("E39").Value = "  +7.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.97"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("E43").Value = "  +4.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.20"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "276.58"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0933"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.81"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.53%  "
